$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) "Sound Storage" paragraph: merge the split "support" / "s" / " up to..."
#    runs back into a single sentence reading "...current setup supports
#    up to 32 individual sounds, placed at 32KB offsets from the start
#    of the ROM."
# ------------------------------------------------------------------
$storagePara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Our design uses an 8Mbit flash memory chip*") {
        $storagePara = $p
        break
    }
}
$null = $storagePara.Range.Find.Execute(
    "Therefore, our current setup supports up to 32 individual sounds, placed at 32KB offsets from the start of the ROM.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Therefore, our current setup supports up to 32 individual sounds, placed at 32KB offsets from the start of the ROM.",
    2)

# ------------------------------------------------------------------
# 2) "Sound Playback" section: move the trailing "To support playback of
#    multiple independent sounds..." sentence so that it becomes part of
#    the first paragraph of that section (appended right after "...multiple
#    individual sounds."), and remove the now-empty trailing paragraphs.
# ------------------------------------------------------------------

# Locate the paragraph that currently holds the standalone sentence, and
# the blank paragraph that immediately precedes it, so both can be removed.
$movedSentencePara = $null
$blankBeforeIndex = $null
$paras = $d.Paragraphs
for ($i = 1; $i -le $paras.Count; $i++) {
    $t = $paras.Item($i).Range.Text
    if ($t -like "To support playback of multiple independent sounds*") {
        $movedSentencePara = $paras.Item($i)
        $blankBeforeIndex = $i - 1
        break
    }
}

$sentenceText = "To support playback of multiple independent sounds, it’s necessary to fetch multiple samples each sample period. The FSM handles this fetching and accumulates the individual sample values in a dedicated register. The I2S interface latches this accumulated sample data each sample period and continuously streams the serialized data to the DAC."

# Delete the sentence paragraph, then the blank paragraph before it (delete
# from the end backwards so indices of earlier paragraphs stay valid).
$movedSentencePara.Range.Delete()
$paras.Item($blankBeforeIndex).Range.Delete()

# Find the "Each sound clip is stored..." paragraph (first paragraph of the
# Sound Playback section) and append the moved sentence to its end, right
# before the paragraph mark.
$playbackIntroPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Each sound clip is stored as a series*") {
        $playbackIntroPara = $p
        break
    }
}
$endPos = $playbackIntroPara.Range.End - 1
$insertionPoint = $d.Range($endPos, $endPos)
$insertionPoint.InsertAfter(" " + $sentenceText)

# ------------------------------------------------------------------
# 3) Update the default run color for the "Normal" and "No Spacing"
#    styles from "automatic" to RGB 00000A.
#    Word color longs are 0xBBGGRR, so RGB 00,00,0A -> 0x0A0000 = 655360.
# ------------------------------------------------------------------
$d.Styles.Item("Normal").Font.Color = 655360
$d.Styles.Item("No Spacing").Font.Color = 655360
